# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# This updates the DAMSLTag (column I) and DialogAct (column J) values for the
# rows whose dialog-act classification changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new DAMSLTag (col I), new DialogAct (col J)
$updates = @(
    @(2,   "b",  "Acknowledge (Backchannel)"),
    @(3,   "aa", "Agree/Accept"),
    @(10,  "sv", "Statement-opinion"),
    @(13,  "sd", "Statement-non-opinion"),
    @(16,  "sd", "Statement-non-opinion"),
    @(21,  "b",  "Acknowledge (Backchannel)"),
    @(22,  "sv", "Statement-opinion"),
    @(27,  "aa", "Agree/Accept"),
    @(32,  "b",  "Acknowledge (Backchannel)"),
    @(33,  "sd", "Statement-non-opinion"),
    @(36,  "sd", "Statement-non-opinion"),
    @(37,  "sd", "Statement-non-opinion"),
    @(48,  "qy", "Yes-No-Question"),
    @(67,  "sd", "Statement-non-opinion"),
    @(71,  "sd", "Statement-non-opinion"),
    @(75,  "qy", "Yes-No-Question"),
    @(81,  "aa", "Agree/Accept"),
    @(91,  "sd", "Statement-non-opinion"),
    @(95,  "sv", "Statement-opinion"),
    @(103, "aa", "Agree/Accept"),
    @(115, "sd", "Statement-non-opinion"),
    @(116, "aa", "Agree/Accept"),
    @(136, "sv", "Statement-opinion"),
    @(137, "b",  "Acknowledge (Backchannel)"),
    @(145, "b",  "Acknowledge (Backchannel)"),
    @(160, "ba", "Appreciation"),
    @(168, "b",  "Acknowledge (Backchannel)"),
    @(174, "sv", "Statement-opinion"),
    @(183, "sd", "Statement-non-opinion"),
    @(184, "sv", "Statement-opinion"),
    @(191, "b",  "Acknowledge (Backchannel)"),
    @(193, "sd", "Statement-non-opinion"),
    @(196, "sv", "Statement-opinion"),
    @(200, "sv", "Statement-opinion"),
    @(201, "sv", "Statement-opinion"),
    @(205, "sv", "Statement-opinion"),
    @(217, "sv", "Statement-opinion"),
    @(222, "sd", "Statement-non-opinion"),
    @(227, "aa", "Agree/Accept"),
    @(232, "sd", "Statement-non-opinion"),
    @(235, "b",  "Acknowledge (Backchannel)"),
    @(236, "sv", "Statement-opinion"),
    @(277, "aa", "Agree/Accept"),
    @(278, "sv", "Statement-opinion"),
    @(287, "sd", "Statement-non-opinion"),
    @(288, "ba", "Appreciation"),
    @(303, "aa", "Agree/Accept"),
    @(321, "sd", "Statement-non-opinion"),
    @(322, "sv", "Statement-opinion"),
    @(323, "sv", "Statement-opinion"),
    @(332, "%",  "Uninterpretable"),
    @(361, "ba", "Appreciation"),
    @(363, "sv", "Statement-opinion"),
    @(365, "sv", "Statement-opinion"),
    @(369, "sv", "Statement-opinion"),
    @(375, "aa", "Agree/Accept"),
    @(380, "sv", "Statement-opinion")
)

foreach ($u in $updates) {
    $row = $u[0]
    $damslTag = $u[1]
    $dialogAct = $u[2]
    $ws.Cells.Item($row, 9).Value = $damslTag
    $ws.Cells.Item($row, 10).Value = $dialogAct
}
